$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 19: 2022-12-20 03:11:04 | Ansh Chawla | ENTRY | Card RESCAN ---
$ws.Range("A19").Value = "2022-12-20 03:11:04"
$ws.Range("B19").Value = 781678351607
$ws.Range("C19").Value = "Ansh Chawla"
$ws.Range("D19").Value = "'19105031"
$ws.Range("E19").Value = "'7696046760"
$ws.Range("F19").Value = "ENTRY"
$ws.Range("G19").Value = "Card RESCAN"

# --- Row 20: 2022-12-20 03:11:36 | Ansh Chawla | ENTRY | Card RESCAN ---
$ws.Range("A20").Value = "2022-12-20 03:11:36"
$ws.Range("B20").Value = 781678351607
$ws.Range("C20").Value = "Ansh Chawla"
$ws.Range("D20").Value = "'19105031"
$ws.Range("E20").Value = "'7696046760"
$ws.Range("F20").Value = "ENTRY"
$ws.Range("G20").Value = "Card RESCAN"

# --- Row 21: 2022-12-20 03:12:09 | Ansh Chawla | EXIT (no Location) ---
$ws.Range("A21").Value = "2022-12-20 03:12:09"
$ws.Range("B21").Value = 781678351607
$ws.Range("C21").Value = "Ansh Chawla"
$ws.Range("D21").Value = "'19105031"
$ws.Range("E21").Value = "'7696046760"
$ws.Range("F21").Value = "EXIT"

# --- Row 22: 2022-12-20 03:12:20 | Ansh Chawla | EXIT | Card RESCAN ---
$ws.Range("A22").Value = "2022-12-20 03:12:20"
$ws.Range("B22").Value = 781678351607
$ws.Range("C22").Value = "Ansh Chawla"
$ws.Range("D22").Value = "'19105031"
$ws.Range("E22").Value = "'7696046760"
$ws.Range("F22").Value = "EXIT"
$ws.Range("G22").Value = "Card RESCAN"
